# Fix sensing_depth column (I) values: the workbook hard-coded 500 for
# every row instead of the correct computed value of 400. Update all
# data rows (2 through 92) in column I from 500 to 400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 9).End(-4162).Row  # xlUp = -4162, column I = 9

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq 500) {
        $cell.Value = 400
    }
}
